$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Master" column (B) used to hold city names (e.g. "Golden State"); it now
# holds the team nickname (e.g. "Warriors") instead. Row order and every other
# column (Start Year, win% figures, etc.) are unchanged -- this only relabels
# the text shown in column B for rows 2..181.
# ---------------------------------------------------------------------------
$teamNames = @(
    "Warriors", "Spurs", "Cavaliers", "Thunder", "Raptors", "Clippers", "Hornets", "Heat", "Celtics", "Hawks",  # rows 2-11
    "Pacers", "Trail Blazers", "Bulls", "Pistons", "Wizards", "Mavericks", "Grizzlies", "Jazz", "Rockets", "Magic",  # rows 12-21
    "Bucks", "Nuggets", "Kings", "Knicks", "Pelicans", "Timberwolves", "Suns", "Nets", "Lakers", "76ers",  # rows 22-31
    "Warriors", "Hawks", "Cavaliers", "Rockets", "Clippers", "Grizzlies", "Spurs", "Trail Blazers", "Bulls", "Mavericks",  # rows 32-41
    "Raptors", "Wizards", "Thunder", "Pelicans", "Bucks", "Suns", "Celtics", "Jazz", "Pacers", "Nets",  # rows 42-51
    "Heat", "Hornets", "Pistons", "Nuggets", "Kings", "Magic", "Lakers", "76ers", "Knicks", "Timberwolves",  # rows 52-61
    "Spurs", "Thunder", "Heat", "Clippers", "Pacers", "Rockets", "Trail Blazers", "Warriors", "Grizzlies", "Suns",  # rows 62-71
    "Mavericks", "Raptors", "Bulls", "Wizards", "Nets", "Hornets", "Timberwolves", "Hawks", "Knicks", "Nuggets",  # rows 72-81
    "Pelicans", "Cavaliers", "Pistons", "Kings", "Lakers", "Celtics", "Jazz", "Magic", "76ers", "Bucks",  # rows 82-91
    "Heat", "Spurs", "Thunder", "Nuggets", "Grizzlies", "Clippers", "Knicks", "Pacers", "Nets", "Warriors",  # rows 92-101
    "Rockets", "Bulls", "Jazz", "Lakers", "Hawks", "Mavericks", "Celtics", "Bucks", "Raptors", "76ers",  # rows 102-111
    "Trail Blazers", "Timberwolves", "Wizards", "Pistons", "Kings", "Pelicans", "Suns", "Cavaliers", "Hornets", "Magic",  # rows 112-121
    "Spurs", "Bulls", "Thunder", "Heat", "Pacers", "Grizzlies", "Lakers", "Hawks", "Celtics", "Clippers",  # rows 122-131
    "Nuggets", "Magic", "76ers", "Knicks", "Rockets", "Mavericks", "Jazz", "Suns", "Bucks", "Trail Blazers",  # rows 132-141
    "Timberwolves", "Pistons", "Raptors", "Warriors", "Kings", "Nets", "Cavaliers", "Pelicans", "Wizards", "Hornets",  # rows 142-151
    "Bulls", "Spurs", "Mavericks", "Heat", "Celtics", "Lakers", "Thunder", "Magic", "Nuggets", "Trail Blazers",  # rows 152-161
    "Grizzlies", "Pelicans", "Hawks", "Rockets", "Knicks", "Suns", "76ers", "Jazz", "Warriors", "Pacers",  # rows 162-171
    "Bucks", "Hornets", "Clippers", "Pistons", "Kings", "Nets", "Wizards", "Raptors", "Cavaliers", "Timberwolves"  # rows 172-181
)

for ($i = 0; $i -lt $teamNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $teamNames[$i]
}

# Restore the scrolled position / active selection recorded in the sheet view
# (best effort -- some window-state attributes are session-local and are
# re-derived by the host on save regardless of what we set here).
$excel.ActiveWindow.ScrollRow = 152
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B182").Select()
